$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new trading day (2021-11-22) as row 14, inheriting the
# "last row" date format that row 13 currently uses.
$ws.Range("A14").NumberFormat = $ws.Range("A13").NumberFormat
$ws.Range("A14").Value = 44522
$ws.Range("B14").Value = 63372.8

# Row 13 is no longer the last row, so it reverts to the regular
# date format used by the other data rows.
$ws.Range("A13").NumberFormat = $ws.Range("A12").NumberFormat
